$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume 1h) are stored as text in this sheet
# (prices use "." as thousands separators, percentages carry padding spaces).
# Force text format before assigning so COM does not silently coerce the
# numeric-looking strings (e.g. "0.998", "42.00", "0.0000263") into Doubles.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "66.608.82"
$ws.Range("E2").Value = "  -3.18%  "
$ws.Range("D3").Value = "3.563.97"
$ws.Range("E3").Value = "  -4.06%  "
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "572.96"
$ws.Range("E5").Value = "  -7.63%  "
$ws.Range("D6").Value = "187.67"
$ws.Range("E6").Value = "  -3.58%  "
$ws.Range("D7").Value = "3.558.79"
$ws.Range("E7").Value = "  -4.14%  "
$ws.Range("D8").Value = "0.616"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "0.673"
$ws.Range("E10").Value = "  -7.46%  "
$ws.Range("E11").Value = "  -7.03%  "
$ws.Range("D12").Value = "55.49"
$ws.Range("E12").Value = "  -7.49%  "
$ws.Range("D13").Value = "0.0000263"
$ws.Range("E13").Value = "  -8.64%  "
$ws.Range("D14").Value = "9.79"
$ws.Range("E14").Value = "  -6.38%  "
$ws.Range("D15").Value = "4.123.40"
$ws.Range("E15").Value = "  -4.35%  "
$ws.Range("D16").Value = "3.560.92"
$ws.Range("E16").Value = "  -4.26%  "
$ws.Range("E17").Value = "  -1.68%  "
$ws.Range("D18").Value = "18.32"
$ws.Range("E18").Value = "  -5.88%  "
$ws.Range("D19").Value = "66.545.77"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").Value = "12.11"
$ws.Range("E20").Value = "  -6.06%  "
$ws.Range("E21").Value = "  -8.06%  "
$ws.Range("D22").Value = "394.11"
$ws.Range("E22").Value = "  -4.16%  "
$ws.Range("D23").Value = "4.19"
$ws.Range("E23").Value = "  -10.66%  "
$ws.Range("D24").Value = "85.36"
$ws.Range("E24").Value = "  -5.42%  "
$ws.Range("D25").Value = "11.28"
$ws.Range("E25").Value = "  -2.55%  "
$ws.Range("D26").Value = "2.92"
$ws.Range("E26").Value = "  -5.64%  "
$ws.Range("D27").Value = "12.41"
$ws.Range("E27").Value = "  -5.21%  "
$ws.Range("D28").Value = "6.05"
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "3.58"
$ws.Range("E29").Value = "  -5.51%  "
$ws.Range("D30").Value = "8.85"
$ws.Range("E30").Value = "  -8.64%  "
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").Value = "30.87"
$ws.Range("E32").Value = "  -6.01%  "
$ws.Range("D33").Value = "635.37"
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("D34").Value = "12.16"
$ws.Range("E34").Value = "  -4.48%  "
$ws.Range("D35").Value = "0.114"
$ws.Range("E35").Value = "  -7.86%  "
$ws.Range("E36").Value = "  -6.33%  "
$ws.Range("D37").Value = "42.00"
$ws.Range("E37").Value = "  -9.95%  "
$ws.Range("D38").Value = "0.402"
$ws.Range("E38").Value = "  -3.40%  "
$ws.Range("E39").Value = "  +0.25%  "
$ws.Range("D40").Value = "0.0₃0751"
$ws.Range("E40").Value = "  -9.34%  "
$ws.Range("D41").Value = "3.136.32"
$ws.Range("E41").Value = "  +6.53%  "
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("E43").Value = "  -0.23%  "
$ws.Range("D44").Value = "2.67"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("D45").Value = "2.94"
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("D46").Value = "0.0412"
$ws.Range("E46").Value = "  -8.07%  "
$ws.Range("E47").Value = "  -6.47%  "
$ws.Range("D48").Value = "3.07"
$ws.Range("E48").Value = "  +0.10%  "
$ws.Range("D49").Value = "139.72"
$ws.Range("E49").Value = "  -4.31%  "
$ws.Range("D50").Value = "8.42"
$ws.Range("E50").Value = "  -11.29%  "
$ws.Range("D51").Value = "2.74"
$ws.Range("E51").Value = "  -1.03%  "
